$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "294.05"
Set-TextValue $ws.Range("E2") "1.27%"
Set-TextValue $ws.Range("D3") "40.10"
Set-TextValue $ws.Range("E3") "1.04%"
Set-TextValue $ws.Range("D4") "5.022"
Set-TextValue $ws.Range("E4") "0.14%"
Set-TextValue $ws.Range("D5") "0.07320"
Set-TextValue $ws.Range("E5") "-0.33%"
Set-TextValue $ws.Range("B6") "GateToken"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D6") "4.316"
Set-TextValue $ws.Range("E6") "0.64%"
Set-TextValue $ws.Range("B7") "FTXToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D7") "1.546"
Set-TextValue $ws.Range("E7") "-0.35%"
Set-TextValue $ws.Range("B8") "MXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9223"
Set-TextValue $ws.Range("E8") "0.81%"
Set-TextValue $ws.Range("B9") "BTSEToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.401"
Set-TextValue $ws.Range("E9") "0.16%"
Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1171"
Set-TextValue $ws.Range("E10") "-1.40%"
Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1783"
Set-TextValue $ws.Range("E11") "2.71%"
Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.08728"
Set-TextValue $ws.Range("E12") "0.15%"
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.04254"
Set-TextValue $ws.Range("E13") "1.79%"
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.1054"
Set-TextValue $ws.Range("E14") "0.39%"
Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001281"
Set-TextValue $ws.Range("E15") "0.16%"
Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005816"
Set-TextValue $ws.Range("E16") "-0.57%"
Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.336"
Set-TextValue $ws.Range("E17") "-1.31%"
Set-TextValue $ws.Range("D18") "0.3293"
Set-TextValue $ws.Range("E18") "-0.09%"
Set-TextValue $ws.Range("D19") "7.948"
Set-TextValue $ws.Range("E19") "5.30%"
Set-TextValue $ws.Range("D20") "0.1383"
Set-TextValue $ws.Range("E20") "2.15%"
Set-TextValue $ws.Range("D21") "0.2809"
Set-TextValue $ws.Range("E21") "2.73%"
Set-TextValue $ws.Range("D22") "0.03948"
Set-TextValue $ws.Range("E22") "2.95%"
Set-TextValue $ws.Range("D23") "0.001271"
Set-TextValue $ws.Range("E23") "-0.47%"
Set-TextValue $ws.Range("D24") "0.003792"
Set-TextValue $ws.Range("E24") "2.70%"
Set-TextValue $ws.Range("E25") "-3.80%"
Set-TextValue $ws.Range("D26") "0.0003731"
Set-TextValue $ws.Range("E26") "0.01%"
Set-TextValue $ws.Range("E38") "1.21%"
Set-TextValue $ws.Range("D39") "0.05085"
Set-TextValue $ws.Range("E39") "1.39%"
Set-TextValue $ws.Range("D40") "0.006079"
Set-TextValue $ws.Range("E40") "18.97%"
Set-TextValue $ws.Range("D41") "0.007761"
Set-TextValue $ws.Range("E41") "1.35%"
Set-TextValue $ws.Range("D42") "0.1294"
Set-TextValue $ws.Range("E42") "2.00%"
Set-TextValue $ws.Range("D43") "0.007407"
Set-TextValue $ws.Range("E43") "-0.01%"
Set-TextValue $ws.Range("D44") "0.007819"
Set-TextValue $ws.Range("E44") "1.80%"
Set-TextValue $ws.Range("D45") "0.2922"
Set-TextValue $ws.Range("E45") "-7.47%"
Set-TextValue $ws.Range("D46") "0.00006189"
Set-TextValue $ws.Range("E46") "-4.92%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.11%"
Set-TextValue $ws.Range("D48") "0.04632"
Set-TextValue $ws.Range("E48") "-81.61%"
Set-TextValue $ws.Range("D49") "0.004210"
Set-TextValue $ws.Range("E49") "0.02%"
Set-TextValue $ws.Range("E50") "0.11%"
Set-TextValue $ws.Range("D51") "0.0002006"
Set-TextValue $ws.Range("E51") "0.11%"
